# The deck currently uses the "Integral" design (clrScheme "Red Violet") as
# its presentation theme (ppt/theme/theme1.xml), while the Notes Master kept
# its own, independent theme copy (ppt/theme/theme2.xml) frozen at the
# default "Office Theme" palette.
#
# This edit switches the presentation's Design back to the stock
# "Office Theme" colour palette, i.e. the 12 theme colours (dk1, lt1, dk2,
# lt2, accent1-6, hlink, folHlink) on the shared ThemeColorScheme are set to
# the standard Office defaults. Because every slide/layout shares the single
# slide master's theme part, touching any one slide's ThemeColorScheme
# updates that shared theme for the whole deck.

$p = $ppt.ActivePresentation

# Office Theme default colour palette, in ThemeColorScheme index order:
# 1 dk1, 2 lt1, 3 dk2, 4 lt2, 5 accent1, 6 accent2, 7 accent3, 8 accent4,
# 9 accent5, 10 accent6, 11 hlink, 12 folHlink.
$officeThemeColors = @(
    0,          # dk1      000000
    16777215,   # lt1      FFFFFF
    6968388,    # dk2      44546A
    15132391,   # lt2      E7E6E6
    13998939,   # accent1  5B9BD5
    3243501,    # accent2  ED7D31
    10855845,   # accent3  A5A5A5
    49407,      # accent4  FFC000
    12874308,   # accent5  4472C4
    4697456,    # accent6  70AD47
    12673797,   # hlink    0563C1
    7491477     # folHlink 954F72
)

$slide = $p.Slides.Item(1)
$themeColors = $slide.ThemeColorScheme

for ($i = 1; $i -le $officeThemeColors.Count; $i++) {
    $themeColors.Item($i).RGB = $officeThemeColors[$i - 1]
}
